# Generate Report for Handoff
# Update Priority ("low" -> "ht") and the "Latest Handoff Datetime" timestamps
# for rows 4-7 on the "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 4-7 -> Priority column E, Latest Handoff Datetime column H
foreach ($r in 4..7) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2017-01-03 09:03:58"
}

# de-de sheet: rows 4-7 -> Priority column E, Latest Handoff Datetime column H
foreach ($r in 4..7) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2017-01-03 09:04:11"
}
